# Auto-generated script applying numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 35714.285
$ws.Range("J75").Value = 35714.285
$ws.Range("L75").Value = 35714.285
$ws.Range("N75").Value = -37586.285
$ws.Range("H78").Value = 35714.285
$ws.Range("J78").Value = 35714.285
$ws.Range("L78").Value = 107142.855
$ws.Range("N78").Value = -116502.855
$ws.Range("H87").Value = 31500
$ws.Range("J87").Value = 31500
$ws.Range("L87").Value = 31500
$ws.Range("N87").Value = -33996
$ws.Range("H90").Value = 31500
$ws.Range("J90").Value = 31500
$ws.Range("L90").Value = 94500
$ws.Range("N90").Value = -106980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11050.492
$ws.Range("I32").Value = 9885.775
$ws.Range("J32").Value = 18556.445
$ws.Range("K32").Value = 9885.775
$ws.Range("L32").Value = 18556.445
$ws.Range("M32").Value = -9598.775
$ws.Range("N32").Value = -19130.445
$ws.Range("H74").Value = 1765.7693
$ws.Range("I74").Value = 1143.3334
$ws.Range("K74").Value = 1143.3334
$ws.Range("M74").Value = -269.3334
$ws.Range("H77").Value = 1765.7693
$ws.Range("I77").Value = 1143.3334
$ws.Range("K77").Value = 5716.666999999999
$ws.Range("M77").Value = -1348.666999999999
$ws.Range("H80").Value = 49097.43
$ws.Range("J80").Value = 49097.43
$ws.Range("L80").Value = 49097.43
$ws.Range("N80").Value = -51093.43
$ws.Range("H83").Value = 49097.43
$ws.Range("J83").Value = 49097.43
$ws.Range("L83").Value = 147292.29
$ws.Range("N83").Value = -157276.29
$ws.Range("H95").Value = 37853.5
$ws.Range("J95").Value = 37853.5
$ws.Range("L95").Value = 37853.5
$ws.Range("N95").Value = -43345.5
$ws.Range("H98").Value = 31115.666
$ws.Range("J98").Value = 31115.666
$ws.Range("L98").Value = 31115.666
$ws.Range("N98").Value = -37105.666
$ws.Range("H101").Value = 48542
$ws.Range("J101").Value = 48542
$ws.Range("L101").Value = 48542
$ws.Range("N101").Value = -55032
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 46664
$ws.Range("J106").Value = 46664
$ws.Range("L106").Value = 46664
$ws.Range("N106").Value = -49188

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 34935.5
$ws.Range("J43").Value = 34935.5
$ws.Range("L43").Value = 34935.5
$ws.Range("N43").Value = -35303.5
$ws.Range("H81").Value = 48311
$ws.Range("J81").Value = 48311
$ws.Range("L81").Value = 48311
$ws.Range("N81").Value = -50307
$ws.Range("H82").Value = 26071
$ws.Range("J82").Value = 33761.332
$ws.Range("L82").Value = 33761.332
$ws.Range("N82").Value = -34483.332
$ws.Range("H84").Value = 48311
$ws.Range("J84").Value = 48311
$ws.Range("L84").Value = 144933
$ws.Range("N84").Value = -154917
$ws.Range("H85").Value = 26071
$ws.Range("J85").Value = 33761.332
$ws.Range("L85").Value = 33761.332
$ws.Range("N85").Value = -36257.332
$ws.Range("H88").Value = 22595
$ws.Range("J88").Value = 22595
$ws.Range("L88").Value = 22595
$ws.Range("N88").Value = -23407
$ws.Range("H91").Value = 22595
$ws.Range("J91").Value = 22595
$ws.Range("L91").Value = 22595
$ws.Range("N91").Value = -25403
$ws.Range("H101").Value = 34935.5
$ws.Range("J101").Value = 34935.5
$ws.Range("L101").Value = 34935.5
$ws.Range("N101").Value = -41425.5
$ws.Range("H106").Value = 46835
$ws.Range("J106").Value = 46835
$ws.Range("L106").Value = 46835
$ws.Range("N106").Value = -49359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 8123.5454
$ws.Range("I102").Value = 9999
$ws.Range("J102").Value = 7936
$ws.Range("K102").Value = 29997
$ws.Range("L102").Value = 23808
$ws.Range("M102").Value = -27563
$ws.Range("N102").Value = -28676
$ws.Range("H131").Value = 4782.5356
$ws.Range("I131").Value = 5940.6113
$ws.Range("J131").Value = 2698
$ws.Range("K131").Value = 17821.8339
$ws.Range("L131").Value = 8094
$ws.Range("M131").Value = -12781.8339
$ws.Range("N131").Value = -18174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 22160.428
$ws.Range("J74").Value = 22160.428
$ws.Range("L74").Value = 22160.428
$ws.Range("N74").Value = -24032.428
$ws.Range("H77").Value = 22160.428
$ws.Range("J77").Value = 22160.428
$ws.Range("L77").Value = 66481.284
$ws.Range("N77").Value = -75841.284
$ws.Range("H96").Value = 35753.668
$ws.Range("J96").Value = 35753.668
$ws.Range("L96").Value = 35753.668
$ws.Range("N96").Value = -41245.668
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H106").Value = 38280
$ws.Range("J106").Value = 38280
$ws.Range("L106").Value = 38280
$ws.Range("N106").Value = -40804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 19515.084
$ws.Range("J81").Value = 19515.084
$ws.Range("L81").Value = 19515.084
$ws.Range("N81").Value = -21511.084
$ws.Range("H84").Value = 19515.084
$ws.Range("J84").Value = 19515.084
$ws.Range("L84").Value = 58545.25199999999
$ws.Range("N84").Value = -68529.25199999999
$ws.Range("H96").Value = 28895.75
$ws.Range("J96").Value = 28895.75
$ws.Range("L96").Value = 28895.75
$ws.Range("N96").Value = -34387.75
$ws.Range("H99").Value = 26507.777
$ws.Range("I99").Value = 17217.5
$ws.Range("J99").Value = 33940
$ws.Range("K99").Value = 17217.5
$ws.Range("L99").Value = 33940
$ws.Range("M99").Value = -14222.5
$ws.Range("N99").Value = -39930
$ws.Range("H102").Value = 48561
$ws.Range("J102").Value = 48561
$ws.Range("L102").Value = 48561
$ws.Range("N102").Value = -55051
$ws.Range("H105").Value = 33161.332
$ws.Range("J105").Value = 33161.332
$ws.Range("L105").Value = 33161.332
$ws.Range("N105").Value = -40149.332
$ws.Range("H106").Value = 36183
$ws.Range("J106").Value = 36183
$ws.Range("L106").Value = 36183
$ws.Range("N106").Value = -38707
$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051
$ws.Range("H129").Value = 36614
$ws.Range("J129").Value = 36614
$ws.Range("L129").Value = 36614
$ws.Range("N129").Value = -46614
$ws.Range("H131").Value = 43318
$ws.Range("J131").Value = 43318
$ws.Range("L131").Value = 43318
$ws.Range("N131").Value = -53398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 23620
$ws.Range("J92").Value = 23620
$ws.Range("L92").Value = 23620
$ws.Range("N92").Value = -28612
$ws.Range("H93").Value = 36014.832
$ws.Range("J93").Value = 36014.832
$ws.Range("L93").Value = 36014.832
$ws.Range("N93").Value = -41006.832
$ws.Range("H97").Value = 34232
$ws.Range("J97").Value = 34232
$ws.Range("L97").Value = 34232
$ws.Range("N97").Value = -36214
$ws.Range("H103").Value = 47876
$ws.Range("J103").Value = 47876
$ws.Range("L103").Value = 47876
$ws.Range("N103").Value = -50220
$ws.Range("H106").Value = 33927.43
$ws.Range("J106").Value = 33927.43
$ws.Range("L106").Value = 33927.43
$ws.Range("N106").Value = -36451.43
$ws.Range("H109").Value = 38377
$ws.Range("J109").Value = 38377
$ws.Range("L109").Value = 38377
$ws.Range("N109").Value = -41151
$ws.Range("H136").Value = 17534.117
$ws.Range("I136").Value = 22774.578
$ws.Range("K136").Value = 68323.734
$ws.Range("M136").Value = -65773.734

Write-Host "Applied all Masamune profit corrections"
